# Applies the "feedback 2 implemented, figures updated" edit:
#  1. Update the cached text of the "datetimeFigureOut" date field
#     (shown on the slide master + every custom layout) from
#     23.03.2022 to 02.05.2022.
#  2. Reposition the "Textfeld 25" textbox (the lone "b" label) on
#     slide 1 to its new offset.

$p = $ppt.ActivePresentation

$oldDate = "23.03.2022"
$newDate = "02.05.2022"

function Update-DateShape($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*") {
            if ($shp.TextFrame.HasText) {
                if ($shp.TextFrame.TextRange.Text -eq $oldDate) {
                    $shp.TextFrame.TextRange.Text = $newDate
                }
            }
        }
    }
}

# 1a. Slide master's own Date Placeholder shape.
Update-DateShape $p.SlideMaster.Shapes

# 1b. Every custom layout's Date Placeholder shape.
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    Update-DateShape $layout.Shapes
}

# 2. Move "Textfeld 25" on slide 1 to its new position.
$slide = $p.Slides.Item(1)
$shapes = $slide.Shapes
for ($i = 1; $i -le $shapes.Count; $i++) {
    $shp = $shapes.Item($i)
    if ($shp.Name -eq "Textfeld 25") {
        $shp.Left = 4465394 / 12700
        $shp.Top = 1595250 / 12700
    }
}
